# Activated alti sensor, added timezone in date
# Adds new rows (164-178) of text entries to the "Translation" sheet,
# mirroring entries that the TouchGFX Designer text editor creates when
# new texts are added to the project (new SingleUseId text ids, plus
# typography/alignment/direction and the actual English ("GB") text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

function Set-TranslationRow($Row, $TextId, $Typography, $Alignment, $Direction, $Text) {
    $ws.Cells.Item($Row, 2).Value = $TextId
    $ws.Cells.Item($Row, 3).Value = $Typography
    $ws.Cells.Item($Row, 4).Value = $Alignment
    $ws.Cells.Item($Row, 5).Value = $Direction
    $ws.Cells.Item($Row, 6).Value = $Text
}

Set-TranslationRow 164 "SingleUseId225" "Default" "Center" "LTR" "Settings"
Set-TranslationRow 165 "SingleUseId226" "Default" "Left"   "LTR" "New Text"
Set-TranslationRow 166 "SingleUseId227" "Default" "Left"   "LTR" "New Text"
Set-TranslationRow 167 "SingleUseId228" "Default" "Left"   "LTR" "New Text"
Set-TranslationRow 168 "SingleUseId229" "Default" "Left"   "LTR" "New Text"
Set-TranslationRow 169 "SingleUseId230" "Default" "Left"   "LTR" "New Text"
Set-TranslationRow 170 "SingleUseId231" "Default" "Left"   "LTR" "New Text"
Set-TranslationRow 171 "SingleUseId232" "Default" "Left"   "LTR" "New Text"
Set-TranslationRow 172 "SingleUseId233" "Default" "Left"   "LTR" "New Text"
Set-TranslationRow 173 "SingleUseId234" "Default" "Left"   "LTR" "New Text"
Set-TranslationRow 174 "SingleUseId235" "Default" "Left"   "LTR" "New Text"
Set-TranslationRow 175 "SingleUseId236" "Small"   "Center" "LTR" "<value>"
Set-TranslationRow 176 "SingleUseId237" "Narrow"  "Center" "LTR" "<value>"
Set-TranslationRow 177 "SingleUseId238" "Default" "Left"   "LTR" "Var name"
Set-TranslationRow 178 "SingleUseId239" "Narrow"  "Left"   "LTR" "Value"
